$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric (e.g. "238.17") must be forced to
# Text format first, otherwise Excel auto-converts them to numbers and
# we lose the exact decimal text (trailing zeros, etc).
$textCells = @("D5", "D7", "D9", "D10", "D11", "D12", "D13", "D14", "D16", "D18", "D21", "D23", "D24", "D27", "D28", "D29", "D31", "D32", "D33", "D34", "D35", "D38", "D39", "D40", "D43", "D45", "D46", "D47")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "41.480.27"
$ws.Range("E2").Value = "  -1.41%  "
$ws.Range("D3").Value = "2.161.07"
$ws.Range("E3").Value = "  -3.17%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "238.17"
$ws.Range("E5").Value = "  -2.25%  "
$ws.Range("E6").Value = "  -3.44%  "
$ws.Range("D7").Value = "71.69"
$ws.Range("E7").Value = "  -2.83%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "0.575"
$ws.Range("E9").Value = "  -4.57%  "
$ws.Range("D10").Value = "39.70"
$ws.Range("E10").Value = "  -6.87%  "
$ws.Range("D11").Value = "0.0903"
$ws.Range("E11").Value = "  -5.71%  "
$ws.Range("D12").Value = "54.05"
$ws.Range("E12").Value = "  -5.01%  "
$ws.Range("D13").Value = "0.0998"
$ws.Range("E13").Value = "  -3.86%  "
$ws.Range("D14").Value = "6.67"
$ws.Range("E14").Value = "  -4.20%  "
$ws.Range("D15").Value = "2.485.77"
$ws.Range("E15").Value = "  -3.08%  "
$ws.Range("D16").Value = "14.07"
$ws.Range("E16").Value = "  -1.96%  "
$ws.Range("D17").Value = "2.157.25"
$ws.Range("E17").Value = "  -2.92%  "
$ws.Range("D18").Value = "0.775"
$ws.Range("E18").Value = "  -7.63%  "
$ws.Range("D19").Value = "41.382.34"
$ws.Range("E19").Value = "  -1.24%  "
$ws.Range("E20").Value = "  -3.21%  "
$ws.Range("D21").Value = "69.63"
$ws.Range("E21").Value = "  -4.55%  "
$ws.Range("E22").Value = "  -7.83%  "
$ws.Range("D23").Value = "9.73"
$ws.Range("E23").Value = "  -13.29%  "
$ws.Range("D24").Value = "226.93"
$ws.Range("E24").Value = "  -1.65%  "
$ws.Range("E25").Value = "  -3.61%  "
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("D27").Value = "10.65"
$ws.Range("E27").Value = "  -6.82%  "
$ws.Range("D28").Value = "3.30"
$ws.Range("E28").Value = "  -8.98%  "
$ws.Range("D29").Value = "2.17"
$ws.Range("E29").Value = "  -4.67%  "
$ws.Range("E30").Value = "  -1.32%  "
$ws.Range("D31").Value = "169.17"
$ws.Range("E31").Value = "  +1.22%  "
$ws.Range("D32").Value = "19.73"
$ws.Range("E32").Value = "  -4.09%  "
$ws.Range("D33").Value = "33.45"
$ws.Range("E33").Value = "  +11.12%  "
$ws.Range("D34").Value = "0.0769"
$ws.Range("E34").Value = "  -4.03%  "
$ws.Range("D35").Value = "5.15"
$ws.Range("E35").Value = "  -8.74%  "
$ws.Range("E36").Value = "  -4.09%  "
$ws.Range("E37").Value = "  -3.84%  "
$ws.Range("D38").Value = "4.27"
$ws.Range("E38").Value = "  -0.91%  "
$ws.Range("D39").Value = "0.0301"
$ws.Range("E39").Value = "  -1.05%  "
$ws.Range("D40").Value = "11.95"
$ws.Range("E40").Value = "  -11.20%  "
$ws.Range("E41").Value = "  -2.93%  "
$ws.Range("E42").Value = "  -6.48%  "
$ws.Range("D43").Value = "58.80"
$ws.Range("E43").Value = "  -9.91%  "
$ws.Range("E44").Value = "  -5.15%  "
$ws.Range("D45").Value = "8.35"
$ws.Range("E45").Value = "  -4.71%  "
$ws.Range("D46").Value = "0.0960"
$ws.Range("E46").Value = "  -4.27%  "
$ws.Range("D47").Value = "95.93"
$ws.Range("E47").Value = "  -8.58%  "
$ws.Range("E48").Value = "  -3.85%  "
$ws.Range("E49").Value = "  -5.42%  "
$ws.Range("E50").Value = "  -8.03%  "
$ws.Range("E51").Value = "  -2.66%  "
